$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for changed rows.
# D-column values are text (e.g. "66.749.59", "1.00") so we force text
# number format before assigning, then reset the style to avoid leaving
# stray direct formatting on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.796.39'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.91%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.501.63'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.51%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.37%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  +4.26%  '

$ws.Range("E9").Value = '  +7.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.34'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.91%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.435'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.11%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.107.99'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.71%  '

$ws.Range("E13").Value = '  -0.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.31%  '

$ws.Range("E15").Value = '  +2.81%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.801.04'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.10%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.507.02'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.89%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.34'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.87%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '396.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.52'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.60%  '

$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("E24").Value = '  +2.30%  '

$ws.Range("E25").Value = '  -0.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.20'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.77%  '

$ws.Range("E27").Value = '  -0.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.38%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.30'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.73%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.47'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.39%  '

$ws.Range("E31").Value = '  -0.19%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.85'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.87%  '

$ws.Range("E33").Value = '  -0.26%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.62'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.88%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '162.82'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.99%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.899'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.45%  '

$ws.Range("E37").Value = '  +0.83%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.86'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.21%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.71'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0746'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.35%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.54'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.18'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.91%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.87'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.84%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0314'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '341.67'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.38%  '

$ws.Range("E48").Value = '  +0.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.859'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.86%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.24%  '

# Rows 43 and 44 swap identity (dogwifhat moves above Maker in the ranking)
# along with updated link / price / volume values. Rank numbers in column A
# are unchanged.
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.64'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.73%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.809.76'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.04%  '
